$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates reflecting refreshed cryptocurrency market data.
# Columns B (Coin) and C (Link) are plain text - safe to assign directly.
# Column D (Price) often looks numeric to Excels type-inference, so we
# force the cell to Text format first to preserve the exact original string
# (avoids turning "0.638" into a float / "0.0000100" into "1E-05", etc.).
# Column E (Volume) values are already padded with spaces/percent signs,
# which Excel does not coerce to numbers, so no extra formatting is required.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.361.24"
$ws.Range("E2").Value = "  -3.39%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.251.21"
$ws.Range("E3").Value = "  -4.24%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.78"
$ws.Range("E5").Value = "  -2.86%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.638"
$ws.Range("E6").Value = "  -4.10%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "71.24"
$ws.Range("E7").Value = "  -2.59%  "

$ws.Range("E8").Value = "  +0.14%  "

$ws.Range("E9").Value = "  -5.41%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.103"
$ws.Range("E10").Value = "  +1.86%  "

$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "58.87"
$ws.Range("E11").Value = "  -2.27%  "

$ws.Range("B12").Value = "Avalanche"
$ws.Range("C12").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "36.91"
$ws.Range("E12").Value = "  +11.76%  "

$ws.Range("E13").Value = "  -2.76%  "

$ws.Range("E14").Value = "  -5.41%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.586.77"
$ws.Range("E15").Value = "  -4.25%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.21"
$ws.Range("E16").Value = "  -7.37%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.878"
$ws.Range("E17").Value = "  -3.18%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.255.52"
$ws.Range("E18").Value = "  -3.87%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.276.75"
$ws.Range("E19").Value = "  -3.42%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000100"
$ws.Range("E20").Value = "  -1.72%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.92"
$ws.Range("E21").Value = "  -4.30%  "

$ws.Range("E22").Value = "  -6.61%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "238.96"
$ws.Range("E23").Value = "  -6.93%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.95"
$ws.Range("E24").Value = "  +2.82%  "

$ws.Range("E25").Value = "  -0.16%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.69"
$ws.Range("E26").Value = "  -1.21%  "

$ws.Range("E27").Value = "  -5.34%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.23"
$ws.Range("E28").Value = "  -3.13%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.16"
$ws.Range("E29").Value = "  -5.05%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "168.85"
$ws.Range("E30").Value = "  -4.91%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.84"
$ws.Range("E31").Value = "  -7.95%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.121"
$ws.Range("E32").Value = "  -6.25%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.128"
$ws.Range("E33").Value = "  -6.05%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0732"
$ws.Range("E34").Value = "  -3.22%  "

$ws.Range("E35").Value = "  -0.54%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.84"
$ws.Range("E36").Value = "  -6.63%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.64"
$ws.Range("E37").Value = "  -4.76%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "22.61"
$ws.Range("E38").Value = "  +18.24%  "

$ws.Range("E39").Value = "  -3.35%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.26"
$ws.Range("E40").Value = "  -4.86%  "

$ws.Range("E41").Value = "  -3.60%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "67.31"
$ws.Range("E42").Value = "  +0.85%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.04"
$ws.Range("E43").Value = "  -0.21%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.93"
$ws.Range("E44").Value = "  -2.25%  "

$ws.Range("E45").Value = "  -8.20%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.190"
$ws.Range("E46").Value = "  -5.54%  "

$ws.Range("E47").Value = "  +0.02%  "

$ws.Range("B48").Value = "Celestia"
$ws.Range("C48").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.61"
$ws.Range("E48").Value = "  +12.22%  "

$ws.Range("B49").Value = "SynthetixNetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.47"
$ws.Range("E49").Value = "  +3.43%  "

$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.40"
$ws.Range("E50").Value = "  -4.02%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.18"
$ws.Range("E51").Value = "  -5.42%  "
